$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 254, pushing the existing rows
# (old 254..350) down to (255..351), matching the new dimension A1:T351.
$ws.Rows(254).Insert()

# Populate the newly inserted row with the new weekly data point.
$ws.Range("A254").Value = 11
$ws.Range("B254").Value = "Vega Monumental Concepción"
$ws.Range("C254").Value = "Bíobío"
$ws.Range("D254").Value = 44825
$ws.Range("E254").Value = 8
$ws.Range("F254").Value = "Fruta"
$ws.Range("G254").Value = 100102
$ws.Range("H254").Value = "Cítricos"
$ws.Range("I254").Value = 100102005
$ws.Range("J254").Value = "Naranja"
$ws.Range("K254").Value = "Lane Late"
$ws.Range("L254").Value = "Primera"
$ws.Range("M254").Value = 200
$ws.Range("N254").Value = 6000
$ws.Range("O254").Value = 6500
$ws.Range("P254").Value = 6250
$ws.Range("Q254").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R254").Value = "Región de O'Higgins"
$ws.Range("S254").Value = 417
$ws.Range("T254").Value = 15
